$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Hola Davivienda", "Hola, gracias por usar el ChatBot de Davivienda >"),
    @("Hola davivienda", "Hola, gracias por usar el ChatBot de Davivienda  >"),
    @("Tarjeta de credito", "Para solicitar una tarjeta de credito, por favor llamar al numero 234  >"),
    @("Tarjeta de crédito", "Para solicitar una tarjeta de credito, por favor llamar al numero 234  >"),
    @("Davivienda", "Hola, gracias por usar el ChatBot de Davivienda  >"),
    @("Me podrias decir por cual empresa fuiste desarrollada", "Claro, por Creasistemas  >"),
    @("Quien te creo?", "Creasistemas  >"),
    @("Quien te ajusto?", "Creasistemas  >"),
    @("Cual empresa te usa?", "Creasistemas  >"),
    @("El banco para el cual estas desarrollada", "Davivienda >"),
    @("Podrias decirme hola", "Claro, hola usuario de Creasistemas es un placer ayudarte >")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Range("B9").Select()
$excel.ActiveWindow.Zoom = 144
